# Apply cell value updates to Sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
# Data-only changes (cached numeric values); no formulas involved.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2928382
$ws.Range("I43").Value = 4050413.5
$ws.Range("J43").Value = 11099.6
$ws.Range("K43").Value = 4050413.5
$ws.Range("L43").Value = 11099.6
$ws.Range("M43").Value = -4050344.5
$ws.Range("N43").Value = -11237.6
$ws.Range("H106").Value = 30514.225
$ws.Range("I106").Value = 33010.543
$ws.Range("K106").Value = 33010.543
$ws.Range("M106").Value = -32379.543
$ws.Range("H112").Value = 215113.55
$ws.Range("J112").Value = 234750.98
$ws.Range("L112").Value = 704252.9400000001
$ws.Range("N112").Value = -706468.9400000001
$ws.Range("H132").Value = 2700.3618
$ws.Range("I132").Value = 2766.5122
$ws.Range("J132").Value = 2248.3333
$ws.Range("K132").Value = 8299.536599999999
$ws.Range("L132").Value = 6744.999899999999
$ws.Range("M132").Value = -5769.536599999999
$ws.Range("N132").Value = -11804.9999
$ws.Range("H138").Value = 530876.6
$ws.Range("I138").Value = 3333
$ws.Range("J138").Value = 629791.0600000001
$ws.Range("K138").Value = 9999
$ws.Range("L138").Value = 1889373.18
$ws.Range("M138").Value = -4859
$ws.Range("N138").Value = -1899653.18

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 21579.6
$ws.Range("I31").Value = 9474.5
$ws.Range("K31").Value = 9474.5
$ws.Range("M31").Value = -9180.5
$ws.Range("H45").Value = 2471.96
$ws.Range("I45").Value = 1865.1
$ws.Range("J45").Value = 4899.4
$ws.Range("K45").Value = 1865.1
$ws.Range("L45").Value = 4899.4
$ws.Range("M45").Value = -1488.1
$ws.Range("N45").Value = -5653.4
$ws.Range("H74").Value = 3083.5117
$ws.Range("I74").Value = 1963.4814
$ws.Range("K74").Value = 1963.4814
$ws.Range("M74").Value = -1089.4814
$ws.Range("H77").Value = 3083.5117
$ws.Range("I77").Value = 1963.4814
$ws.Range("K77").Value = 9817.406999999999
$ws.Range("M77").Value = -5449.406999999999
$ws.Range("H98").Value = 100000
$ws.Range("J98").Value = 100000
$ws.Range("L98").Value = 100000
$ws.Range("N98").Value = -105990

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H82").Value = 12877.091
$ws.Range("J82").Value = 39712
$ws.Range("L82").Value = 39712
$ws.Range("N82").Value = -40478
$ws.Range("H85").Value = 12877.091
$ws.Range("J85").Value = 39712
$ws.Range("L85").Value = 39712
$ws.Range("N85").Value = -42364
$ws.Range("H105").Value = 2599.875
$ws.Range("I105").Value = 879.75
$ws.Range("K105").Value = 879.75
$ws.Range("M105").Value = 867.25
$ws.Range("H134").Value = 5945.3076
$ws.Range("I134").Value = 6429
$ws.Range("K134").Value = 19287
$ws.Range("M134").Value = -16752

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H31").Value = 1838.4
$ws.Range("I31").Value = 1323.25
$ws.Range("K31").Value = 1323.25
$ws.Range("M31").Value = -1028.25
$ws.Range("H34").Value = 1838.4
$ws.Range("I34").Value = 1323.25
$ws.Range("K34").Value = 1323.25
$ws.Range("M34").Value = -1121.25
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H58").Value = 4660.6113
$ws.Range("I58").Value = 5234
$ws.Range("K58").Value = 5234
$ws.Range("M58").Value = -5031
$ws.Range("H95").Value = 65406
$ws.Range("J95").Value = 65406
$ws.Range("L95").Value = 65406
$ws.Range("N95").Value = -70898
$ws.Range("H132").Value = 6615.154
$ws.Range("I132").Value = 6499.5
$ws.Range("J132").Value = 6636.1816
$ws.Range("K132").Value = 19498.5
$ws.Range("L132").Value = 19908.5448
$ws.Range("M132").Value = -16968.5
$ws.Range("N132").Value = -24968.5448
$ws.Range("H134").Value = 5072.7896
$ws.Range("I134").Value = 4185.478
$ws.Range("K134").Value = 12556.434
$ws.Range("M134").Value = -10021.434
$ws.Range("H136").Value = 4660.6113
$ws.Range("I136").Value = 5234
$ws.Range("K136").Value = 15702
$ws.Range("M136").Value = -13152

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2131.375
$ws.Range("J68").Value = 2283.6667
$ws.Range("L68").Value = 6851.000100000001
$ws.Range("N68").Value = -8473.000100000001
$ws.Range("H71").Value = 2131.375
$ws.Range("J71").Value = 2283.6667
$ws.Range("L71").Value = 20553.0003
$ws.Range("N71").Value = -28665.0003
$ws.Range("H129").Value = 3319.7932
$ws.Range("J129").Value = 4503
$ws.Range("L129").Value = 13509
$ws.Range("N129").Value = -23509

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()
$ws.Range("H93").Value = 87919
$ws.Range("I93").Value = 87888
$ws.Range("J93").Value = 87925.2
$ws.Range("K93").Value = 87888
$ws.Range("L93").Value = 87925.2
$ws.Range("M93").Value = -86016
$ws.Range("N93").Value = -91669.2
$ws.Range("H97").Value = 1627.7059
$ws.Range("I97").Value = 917.7857
$ws.Range("J97").Value = 4940.6665
$ws.Range("K97").Value = 917.7857
$ws.Range("L97").Value = 4940.6665
$ws.Range("M97").Value = -421.7857
$ws.Range("N97").Value = -5932.6665
$ws.Range("H99").Value = 19179.059
$ws.Range("I99").Value = 1681.75
$ws.Range("J99").Value = 34732.223
$ws.Range("K99").Value = 1681.75
$ws.Range("L99").Value = 34732.223
$ws.Range("M99").Value = 564.25
$ws.Range("N99").Value = -39224.223
$ws.Range("H132").Value = 1961.88
$ws.Range("I132").Value = 1968.3478
$ws.Range("J132").Value = 1887.5
$ws.Range("K132").Value = 5905.0434
$ws.Range("L132").Value = 5662.5
$ws.Range("M132").Value = -3375.0434
$ws.Range("N132").Value = -10722.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6000.148
$ws.Range("I46").Value = 3773.3333
$ws.Range("K46").Value = 3773.3333
$ws.Range("M46").Value = -3585.3333
$ws.Range("H55").Value = 709.4194
$ws.Range("I55").Value = 526.3
$ws.Range("J55").Value = 1042.3636
$ws.Range("K55").Value = 526.3
$ws.Range("L55").Value = 1042.3636
$ws.Range("M55").Value = -353.3
$ws.Range("N55").Value = -1388.3636
$ws.Range("H100").Value = 3682.111
$ws.Range("I100").Value = 3302.4443
$ws.Range("J100").Value = 4441.4443
$ws.Range("K100").Value = 3302.4443
$ws.Range("L100").Value = 4441.4443
$ws.Range("M100").Value = -2761.4443
$ws.Range("N100").Value = -5523.4443
$ws.Range("H136").Value = 3304.5454
$ws.Range("I136").Value = 3095.9443
$ws.Range("J136").Value = 4243.25
$ws.Range("K136").Value = 9287.832900000001
$ws.Range("L136").Value = 12729.75
$ws.Range("M136").Value = -6737.832900000001
$ws.Range("N136").Value = -17829.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1609.275
$ws.Range("I136").Value = 1520.6666
$ws.Range("J136").Value = 1875.1
$ws.Range("K136").Value = 4561.9998
$ws.Range("L136").Value = 5625.299999999999
$ws.Range("M136").Value = -2011.9998
$ws.Range("N136").Value = -10725.3
$ws.Range("H138").Value = 99623
$ws.Range("J138").Value = 99623
$ws.Range("L138").Value = 99623
$ws.Range("N138").Value = -109903
